# Junction_Flooding_422 data refresh: new simulation run (custom accuracy + 1000 data points)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old trailing data row (row 6) so the sheet shrinks back to A1:AH5
$ws.Rows.Item(6).Delete()

# Replace rows 2-5 (Time + J1..AH32) with the refreshed simulation output
$data = @(
    @(45103.50694444445,20.178,13.652,4.068,42.752,34.691,15.879,51.007,24.432,10.251,15.607,16.868,17.597,5.069,15.79,22.094,13.41,3.498,2.451,232.996,43.923,14.575,29.082,14.996,3.13,25.159,12.874,11.654,13.651,17.331,3.457,45.214,8.105,18.222),
    @(45103.51388888889,0.842,0.294,1.291,1.753,1.011,0.765,10.351,1.163,0.394,0.237,0.68,0.529,0.289,0.752,1.083,1.012,1.436,0.629,4.313,2.762,0.694,1.697,0.716,0.6889999999999999,4.238,0.613,0.757,0.828,0.474,1.266,10.074,0.237,0.888),
    @(45103.52083333334,12.39,9.098000000000001,1.19,26.997,21.972,9.83,34.244,15.125,6.544,9.73,10.849,11.334,3.136,9.775,13.747,8.411,1.024,0.709,141.419,27.144,9.023,18.038,9.557,1.69,16.63,7.97,7.153,8.388999999999999,11.251,0.773,30.602,5.018,11.28),
    @(45103.52777777778,24.41,18.2,1.38,53.19,43.68,19.28,72.16,29.67,13.17,19.57,21.37,22.49,6.16,19.17,27.28,16.14,0.89,0.96,284.54,53.58,17.7,36.05,18.96,2.83,35.48,15.63,13.85,16.27,22.36,0.5600000000000001,65.52,9.99,22.13)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowVals[$j]
    }
}

# Column width tweaks (several columns widened from 7 to 8, column T from 8 to 9).
# The COM layer reports/stores ColumnWidth with Excel's usual +5/6 char padding,
# so subtract that offset to land on the exact raw OOXML "width" the diff wants.
$widthOffset = 0.8333333333333334
$widths = @{
    1=21; 2=8; 3=8; 4=7; 5=8; 6=8; 7=8; 8=8; 9=8; 10=8; 11=8; 12=8; 13=8;
    14=7; 15=7; 16=8; 17=7; 18=7; 19=7; 20=9; 21=8; 22=8; 23=8; 24=8; 25=7;
    26=8; 27=8; 28=8; 29=8; 30=8; 31=7; 32=8; 33=7; 34=8
}
foreach ($colIdx in $widths.Keys) {
    $ws.Columns.Item($colIdx).ColumnWidth = $widths[$colIdx] - $widthOffset
}
